$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N ("Late" and everything to its right shifts
# one column over: old N->O, old O->P, old P->Q).
$ws.Columns("N:N").Insert()

# Bring this sheet to the front (it becomes the active tab / active sheet)
# and leave the selection where the author left it.
$ws.Activate()
$ws.Range("J19").Select()
